# Update the pl_mw.xlsx "Sheet1" data with new values (case with 380 kV done).
$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Sheet1")

# Columns that change, in order: B, C, E, F, G, H, M, N
$cols = @("B","C","E","F","G","H","M","N")

# Row data: row number followed by the 8 new values (for columns B,C,E,F,G,H,M,N)
$rowsData = @(
    @(2,  0.8007386384252868, 0.1891880871935996, 0.114491997644862,  0.4443680307746121, 0.558292718289735,  0.6226879084319137, 0.3383916136339522, 1.008337260620687),
    @(3,  0.705321384394324,  0.1663163533809495, 0.1086070773511167, 0.3878228170618172, 0.5392220746471423, 0.6195508052312135, 0.3020394082274862, 1.02523291381426),
    @(4,  0.6468390329911529, 0.1522137044287604, 0.1050823239116667, 0.3531389305169483, 0.5280821386504897, 0.6180934530629685, 0.2798358380511985, 1.036132427887289),
    @(5,  0.6230328464602337, 0.1464517051098824, 0.1036679966319767, 0.3390132514313251, 0.5236844113795627, 0.6176169696966127, 0.2708165993970866, 1.040705951671102),
    @(6,  0.6190813981837096, 0.1454940161764853, 0.1034344726387246, 0.336668177824194,  0.5229627031712027, 0.6175449262346291, 0.2693206930375709, 1.041473339673453),
    @(7,  0.6465178702493972, 0.1521360572302513, 0.1050631608871058, 0.3529483938344953, 0.5280222566167367, 0.6180865523554218, 0.2797140850732305, 1.036193574390739),
    @(8,  0.7678168855517242, 0.1813141314118241, 0.1124443367185819, 0.4248636149813478, 0.5515981529282641, 0.6215086256681985, 0.325832822422484,  1.014053654698233),
    @(9,  1.006539895636706,  0.2380707708591956, 0.1276337629410662, 0.5661985755041457, 0.6024120485364364, 0.6319630491316843, 0.4172250574418399, 0.9748181828079359),
    @(10, 1.182509939393753,  0.2795052098613553, 0.139246998362502,  0.6702781546542269, 0.6426295197639718, 0.6419610150192057, 0.4849958932013152, 0.9485550494946358),
    @(11, 1.262703134319452,  0.2983012138844003, 0.1446325694007911, 0.7176906081379002, 0.6615718885135209, 0.6470201778217017, 0.5159722179179482, 0.937167008961965),
    @(12, 1.293091599671925,  0.3054114195198281, 0.1466870022064057, 0.7356546913071611, 0.668839489023668,  0.6490100342238634, 0.5277239531569933, 0.9329353735234696),
    @(13, 1.286545957483611,  0.3038804376906512, 0.1462438703322348, 0.7317853510981394, 0.6672700532606939, 0.6485781799414951, 0.5251920342497982, 0.9338431339649853),
    @(14, 1.265202787563624,  0.2988863228835044, 0.1448012859492565, 0.7191683204515869, 0.6621678953062542, 0.6471823968974491, 0.5169386028157845, 0.9368172489778379),
    @(15, 1.25213223231691,   0.295826319845105,  0.1439196273659675, 0.7114413442032514, 0.6590550313115102, 0.6463371015857433, 0.511885976137421,  0.9386495094504799),
    @(16, 1.177272047903898,  0.2782757984098225, 0.1388971274967759, 0.6671810134426437, 0.6414047404888095, 0.6416407192194242, 0.4829745213400969, 0.9493105674624447),
    @(17, 1.131384794565975,  0.2674957215769496, 0.1358424712818831, 0.6400460337215605, 0.6307436712134233, 0.6388909369874796, 0.4652763134530176, 0.9559942554289282),
    @(18, 1.105005145719133,  0.2612903341232879, 0.1340951594357804, 0.6244449056556647, 0.6246725392698806, 0.6373574016236603, 0.455110579487723,  0.9598911478369292),
    @(19, 1.096075754056642,  0.2591884380590557, 0.1335051980485034, 0.6191636801734006, 0.6226273619893163, 0.6368464109640115, 0.4516709919950301, 0.9612195940023414),
    @(20, 1.136268172436644,  0.2686437927044381, 0.1361666443979885, 0.642933953830422,  0.6318722536715029, 0.6391786771564227, 0.4671588850869028, 0.9552773177629827),
    @(21, 1.27147121955619,   0.3003534159857395, 0.1452245977471236, 0.7228739723492197, 0.663663945362913,  0.6475903576042867, 0.5193622432193337, 0.9359414843906668),
    @(22, 1.359957078262823,  0.321034233836798,  0.1512322517825169, 0.7751780083420101, 0.6849933715668328, 0.6535198331825427, 0.5536067362528314, 0.9237753501994979),
    @(23, 1.312719132322059,  0.3100003995025986, 0.1480177330481709, 0.7472568307830727, 0.6735584868393687, 0.6503154452084345, 0.5353180463466174, 0.9302254330506621),
    @(24, 1.134060392525953,  0.2681247739889159, 0.1360200581859061, 0.6416283278902171, 0.6313618408086938, 0.6390484423730527, 0.4663077458350102, 0.9556012760055044),
    @(25, 0.9418613730704806, 0.2227642913527177, 0.123446258690258,  0.5279251897347166, 0.5881656311262873, 0.6287304768997899, 0.3923944686134604, 0.9849833660033251)
)

foreach ($rowEntry in $rowsData) {
    $rowNum = $rowEntry[0]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $col = $cols[$i]
        $value = $rowEntry[$i + 1]
        $ws.Range("$col$rowNum").Value = $value
    }
}

$wb.Save()
